# -----------------------------------------------------------------------
# Adds three new worksheets (plausibility, grammaticallity, insert_reflexive)
# with example-sentence comparison tables, highlights a new note row on
# List1, and updates the active-sheet/selection view state.
# -----------------------------------------------------------------------

function Set-Edge($rng, $edge, $style) {
    if (-not $style) { return }
    $b = $rng.Borders.Item($edge)
    $b.LineStyle = 1
    if ($style -eq 'medium') {
        $b.Weight = -4138
    } else {
        $b.Weight = 2
    }
}

# xlEdgeLeft=7, xlEdgeTop=8, xlEdgeBottom=9, xlEdgeRight=10
function Set-Box($ws, $addr, $left, $right, $top, $bottom) {
    $rng = $ws.Range($addr)
    Set-Edge $rng 7  $left
    Set-Edge $rng 10 $right
    Set-Edge $rng 8  $top
    Set-Edge $rng 9  $bottom
}

function Fill-Table($ws, $data) {
    foreach ($addr in $data.Keys) {
        $ws.Range($addr).Value = $data[$addr]
    }
}

$wb = $excel.ActiveWorkbook
$list1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# List1: new note row with yellow highlight
# ---------------------------------------------------------------------
$list1.Cells.Item(32, 1).Value = "dalsi slovesa, ktera lze pridat"
$list1.Cells.Item(32, 1).Interior.Color = 65535
$list1.Range("D26").Select()

# ---------------------------------------------------------------------
# Sheet: plausibility
# ---------------------------------------------------------------------
$s2 = $wb.Worksheets.Add([System.Type]::Missing, $list1)
$s2.Name = "plausibility"

Fill-Table $s2 @{
    "C1" = "plausible use";   "D1" = "implausible use"
    "A2" = "high-";           "B2" = "reflexive use";     "C2" = "dědeček se probudil ze snu";  "D2" = "budík se probudil ze snu"
    "A3" = "reflexivity bias";"B3" = "non-reflexive use"; "C3" = "budík probudil dědečka ze snu";"D3" = "dědeček probudil budík ze snu"
    "A4" = "low-";            "B4" = "reflexive use";     "C4" = "Petr se zapsal  do sešitu";   "D4" = "Datum se zapsal do sešitu"
    "A5" = "reflexivity bias";"B5" = "non-reflexive use"; "C5" = "petr zapsal datum do sešitu"; "D5" = "datum zapsal petra do sešitu"
}

$s2.Columns.Item(1).ColumnWidth = 12.66
$s2.Columns.Item(2).ColumnWidth = 15.11
$s2.Columns.Item(3).ColumnWidth = 25.55
$s2.Columns.Item(4).ColumnWidth = 25.55

Set-Box $s2 "A1" 'thin'   'thin'   'thin'   $null
Set-Box $s2 "B1" 'thin'   'thin'   'thin'   'thin'
Set-Box $s2 "C1" 'thin'   'thin'   'thin'   'thin'
Set-Box $s2 "D1" 'thin'   'thin'   'thin'   'thin'

Set-Box $s2 "A2" 'medium' 'medium' 'medium' $null
Set-Box $s2 "B2" $null    'thin'   $null    'thin'
Set-Box $s2 "C2" 'thin'   'thin'   $null    'thin'
Set-Box $s2 "D2" 'thin'   'thin'   $null    'thin'

Set-Box $s2 "A3" 'medium' 'medium' $null    'medium'
Set-Box $s2 "B3" $null    'thin'   'thin'   'thin'
Set-Box $s2 "C3" 'thin'   'thin'   'thin'   'thin'
Set-Box $s2 "D3" 'thin'   'thin'   'thin'   'thin'

Set-Box $s2 "A4" 'medium' 'medium' 'medium' $null
Set-Box $s2 "B4" $null    'thin'   'thin'   'thin'
Set-Box $s2 "C4" 'thin'   'thin'   'thin'   'thin'
Set-Box $s2 "D4" 'thin'   'thin'   'thin'   'thin'

Set-Box $s2 "A5" 'medium' 'medium' $null    'medium'
Set-Box $s2 "B5" $null    'thin'   'thin'   'thin'
Set-Box $s2 "C5" 'thin'   'thin'   'thin'   'thin'
Set-Box $s2 "D5" 'thin'   'thin'   'thin'   'thin'

$s2.Range("D19").Select()

# ---------------------------------------------------------------------
# Sheet: grammaticallity
# ---------------------------------------------------------------------
$s3 = $wb.Worksheets.Add([System.Type]::Missing, $s2)
$s3.Name = "grammaticallity"

Fill-Table $s3 @{
    "C1" = "grammatical use"; "D1" = "non-grammatical use"
    "A2" = "high-";           "B2" = "intransitive use";  "C2" = "dědeček se probudil";          "D2" = "dědeček probudil"
    "A3" = "reflexivity bias";"B3" = "transitive use";    "C3" = "budík probudil dědečka";       "D3" = "budík se probudil dědečka"
    "A4" = "low-";            "B4" = "reflexive use";     "C4" = "petr se zapsal do sešitu";     "D4" = "petr zapsal do sešitu"
    "A5" = "reflexivity bias";"B5" = "non-reflexive use"; "C5" = "petr zapsal datum do sešitu";  "D5" = "petr se zapsal datum do sešitu"
}

$s3.Columns.Item(1).ColumnWidth = 12.66
$s3.Columns.Item(2).ColumnWidth = 15.11
$s3.Columns.Item(3).ColumnWidth = 25.55
$s3.Columns.Item(4).ColumnWidth = 25.55

Set-Box $s3 "A1" 'thin'   'thin'   'thin'   $null
Set-Box $s3 "B1" 'thin'   'thin'   'thin'   'thin'
Set-Box $s3 "C1" 'thin'   'thin'   'thin'   'thin'
Set-Box $s3 "D1" 'thin'   'thin'   'thin'   'thin'

Set-Box $s3 "A2" 'medium' 'medium' 'medium' $null
Set-Box $s3 "B2" $null    'thin'   $null    'thin'
Set-Box $s3 "C2" 'thin'   'thin'   $null    'thin'
Set-Box $s3 "D2" 'thin'   'thin'   $null    'thin'

Set-Box $s3 "A3" 'medium' 'medium' $null    'medium'
Set-Box $s3 "B3" $null    'thin'   'thin'   'thin'
Set-Box $s3 "C3" 'thin'   'thin'   'thin'   'thin'
Set-Box $s3 "D3" 'thin'   'thin'   'thin'   'thin'

Set-Box $s3 "A4" 'medium' 'medium' 'medium' $null
Set-Box $s3 "B4" $null    'thin'   'thin'   'thin'
Set-Box $s3 "C4" 'thin'   'thin'   'thin'   'thin'
Set-Box $s3 "D4" 'thin'   'thin'   'thin'   'thin'

Set-Box $s3 "A5" 'medium' 'medium' $null    'medium'
Set-Box $s3 "B5" $null    'thin'   'thin'   'thin'
Set-Box $s3 "C5" 'thin'   'thin'   'thin'   'thin'
Set-Box $s3 "D5" 'thin'   'thin'   'thin'   'thin'

$s3.PageSetup.PaperSize = 9
$s3.PageSetup.Orientation = 1

$s3.Range("A1:D5").Select()

# ---------------------------------------------------------------------
# Sheet: insert_reflexive
# ---------------------------------------------------------------------
$s4 = $wb.Worksheets.Add([System.Type]::Missing, $s3)
$s4.Name = "insert_reflexive"

Fill-Table $s4 @{
    "C1" = "example sentence"
    "A2" = "high-";           "B2" = "reflexive needed";      "C2" = "dědeček _ probudil _ ze _ snu"
    "A3" = "reflexivity bias";"B3" = "reflexive not needed";  "C3" = "budik _ probudil _ dědečka _ ze _ snu"
    "A4" = "low-";            "B4" = "reflexive needed";      "C4" = "petr _ zapsal _ do _ sešitu"
    "A5" = "reflexivity bias";"B5" = "reflexive not needed";  "C5" = "petr _ zapsal _ datum _ do _ sešitu"
}

$s4.Columns.Item(1).ColumnWidth = 12.66
$s4.Columns.Item(2).ColumnWidth = 17.88
$s4.Columns.Item(3).ColumnWidth = 31.44
$s4.Columns.Item(4).ColumnWidth = 25.77

Set-Box $s4 "A1" 'thin'   'thin'   'thin'   $null
Set-Box $s4 "B1" 'thin'   'thin'   'thin'   'thin'
Set-Box $s4 "C1" 'thin'   'thin'   'thin'   'thin'

Set-Box $s4 "A2" 'medium' 'medium' 'medium' $null
Set-Box $s4 "B2" $null    'thin'   $null    'thin'
Set-Box $s4 "C2" 'thin'   'thin'   $null    'thin'

Set-Box $s4 "A3" 'medium' 'medium' $null    'medium'
Set-Box $s4 "B3" $null    'thin'   'thin'   'thin'
Set-Box $s4 "C3" 'thin'   'thin'   'thin'   'thin'

Set-Box $s4 "A4" 'medium' 'medium' 'medium' $null
Set-Box $s4 "B4" $null    'thin'   'thin'   'thin'
Set-Box $s4 "C4" 'thin'   'thin'   'thin'   'thin'

Set-Box $s4 "A5" 'medium' 'medium' $null    'medium'
Set-Box $s4 "B5" $null    'thin'   'thin'   'thin'
Set-Box $s4 "C5" 'thin'   'thin'   'thin'   'thin'

# Column D has no data but the original workbook still touches D1:D5 so the
# sheet dimension extends to D5 and the (blank) cells are materialised.
$s4.Range("D1:D5").Borders.Item(7).LineStyle = -4142

$s4.Range("C15").Select()
$s4.Activate()
